$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,12

$arr[0,0] = 0.4097860205177142
$arr[0,1] = 0.1508574413176973
$arr[0,2] = 0.6203738845577647
$arr[0,3] = 0.2340256789016379
$arr[0,4] = 0
$arr[0,5] = 2.36627318573656
$arr[0,6] = 1.918577420600258
$arr[0,7] = 1.68155054552291
$arr[0,8] = 0.1077126063114058
$arr[0,9] = 0.5563252529647684
$arr[0,10] = 0
$arr[0,11] = 0.3860665536174892

$arr[1,0] = 0.3869794967756945
$arr[1,1] = 0.1448246039694112
$arr[1,2] = 0.6151783919307405
$arr[1,3] = 0.2325711805172332
$arr[1,4] = 0
$arr[1,5] = 2.343667298821401
$arr[1,6] = 1.912888664632703
$arr[1,7] = 1.672644897987176
$arr[1,8] = 0.1073769678851022
$arr[1,9] = 0.5280216225579863
$arr[1,10] = 0
$arr[1,11] = 0.3771023688069874

$arr[2,0] = 0.3732552226685755
$arr[2,1] = 0.141219158328056
$arr[2,2] = 0.6122727078315364
$arr[2,3] = 0.2317792434477326
$arr[2,4] = 0
$arr[2,5] = 2.330684751943863
$arr[2,6] = 1.910001162740258
$arr[2,7] = 1.667718394629865
$arr[2,8] = 0.1072123384018617
$arr[2,9] = 0.5110197940196031
$arr[2,10] = 0
$arr[2,11] = 0.3718284715613791

$arr[3,0] = 0.3677327502520598
$arr[3,1] = 0.139774697581359
$arr[3,2] = 0.6111601645680764
$arr[3,3] = 0.2314819642676014
$arr[3,4] = 0
$arr[3,5] = 2.325619608601144
$arr[3,6] = 1.90897668212628
$arr[3,7] = 1.665846844577281
$arr[3,8] = 0.1071556774615772
$arr[3,9] = 0.5041861982276288
$arr[3,10] = 0
$arr[3,11] = 0.3697372199704034

$arr[4,0] = 0.3668199962817482
$arr[4,1] = 0.1395363419054831
$arr[4,2] = 0.6109797500559893
$arr[4,3] = 0.2314341382483853
$arr[4,4] = 0
$arr[4,5] = 2.3247921458794
$arr[4,6] = 1.90881575936524
$arr[4,7] = 1.665544287262961
$arr[4,8] = 0.107146898809301
$arr[4,9] = 0.5030572116770173
$arr[4,10] = 0
$arr[4,11] = 0.3693934675070949

$arr[5,0] = 0.3731804599516408
$arr[5,1] = 0.1411995775150956
$arr[5,2] = 0.6122574139399575
$arr[5,3] = 0.2317751312087175
$arr[5,4] = 0
$arr[5,5] = 2.33061552968627
$arr[5,6] = 1.909986730044693
$arr[5,7] = 1.667692603636659
$arr[5,8] = 0.107211532030103
$arr[5,9] = 0.5109272498839914
$arr[5,10] = 0
$arr[5,11] = 0.3718000337147913

$arr[6,0] = 0.4018644809007412
$arr[6,1] = 0.1487567976642623
$arr[6,2] = 0.6185234780629969
$arr[6,3] = 0.2335031838922319
$arr[6,4] = 0
$arr[6,5] = 2.358292166166649
$arr[6,6] = 1.916490255388055
$arr[6,7] = 1.678367373460496
$arr[6,8] = 0.1075882758524038
$arr[6,9] = 0.5464879974888959
$arr[6,10] = 0
$arr[6,11] = 0.3829279498500924

$arr[7,0] = 0.4603262308855847
$arr[7,1] = 0.1643628578461005
$arr[7,2] = 0.6330672387528011
$arr[7,3] = 0.2376942717069355
$arr[7,4] = 0
$arr[7,5] = 2.419708319390793
$arr[7,6] = 1.934051327117288
$arr[7,7] = 1.703607066646597
$arr[7,8] = 0.1086560096811908
$arr[7,9] = 0.6192137078429312
$arr[7,10] = 0
$arr[7,11] = 0.4065762814553722

$arr[8,0] = 0.5046304243093118
$arr[8,1] = 0.1763135471478847
$arr[8,2] = 0.645129566068988
$arr[8,3] = 0.2412632340024388
$arr[8,4] = 0
$arr[8,5] = 2.469220564738833
$arr[8,6] = 1.949893128302875
$arr[8,7] = 1.72479215044217
$arr[8,8] = 0.1096412642205138
$arr[8,9] = 0.674478680470969
$arr[8,10] = 0
$arr[8,11] = 0.4250673935122435

$arr[9,0] = 0.525080375615147
$arr[9,1] = 0.1818568208604177
$arr[9,2] = 0.6509165672180188
$arr[9,3] = 0.2429933822659365
$arr[9,4] = 0
$arr[9,5] = 2.492706316598998
$arr[9,6] = 1.957740439210909
$arr[9,7] = 1.735007073429756
$arr[9,8] = 0.1101331562459791
$arr[9,9] = 0.7000208891273587
$arr[9,10] = 0
$arr[9,11] = 0.4337228402673006

$arr[10,0] = 0.5328667552636261
$arr[10,1] = 0.1839713552601552
$arr[10,2] = 0.6531510676435062
$arr[10,3] = 0.2436638778091691
$arr[10,4] = 0
$arr[10,5] = 2.501738667636857
$arr[10,6] = 1.960804271905914
$arr[10,7] = 1.738958507642934
$arr[10,8] = 0.1103257084417564
$arr[10,9] = 0.7097509294144686
$arr[10,10] = 0
$arr[10,11] = 0.4370355052211892

$arr[11,0] = 0.5311879337166374
$arr[11,1] = 0.1835152659159576
$arr[11,2] = 0.6526679119235439
$arr[11,3] = 0.2435187930351432
$arr[11,4] = 0
$arr[11,5] = 2.499787209100418
$arr[11,6] = 1.960140317612769
$arr[11,7] = 1.738103789347022
$arr[11,8] = 0.1102839594228087
$arr[11,9] = 0.7076528226932339
$arr[11,10] = 0
$arr[11,11] = 0.4363205055787915

$arr[12,0] = 0.5257201160245302
$arr[12,1] = 0.1820304755929101
$arr[12,2] = 0.6510995374616471
$arr[12,3] = 0.2430482372058549
$arr[12,4] = 0
$arr[12,5] = 2.493446629267254
$arr[12,6] = 1.957990653747629
$arr[12,7] = 1.735330490468641
$arr[12,8] = 0.1101488717014192
$arr[12,9] = 0.7008202279600084
$arr[12,10] = 0
$arr[12,11] = 0.4339946731315294

$arr[13,0] = 0.5223764420962027
$arr[13,1] = 0.1811230079790676
$arr[13,2] = 0.6501444736605606
$arr[13,3] = 0.2427620039446694
$arr[13,4] = 0
$arr[13,5] = 2.489580930168984
$arr[13,6] = 1.956685935584744
$arr[13,7] = 1.733642613959233
$arr[13,8] = 0.1100669448574081
$arr[13,9] = 0.6966425876318851
$arr[13,10] = 0
$arr[13,11] = 0.4325745957857805

$arr[14,0] = 0.5032999143734571
$arr[14,1] = 0.1759534356480685
$arr[14,2] = 0.6447574026468601
$arr[14,3] = 0.241152309612481
$arr[14,4] = 0
$arr[14,5] = 2.467705116810066
$arr[14,6] = 1.94939318938043
$arr[14,7] = 1.724136222559025
$arr[14,8] = 0.1096099971743811
$arr[14,9] = 0.6728175213793861
$arr[14,10] = 0
$arr[14,11] = 0.4245066431591056

$arr[15,0] = 0.4916727577701749
$arr[15,1] = 0.1728094731599299
$arr[15,2] = 0.6415293768519348
$arr[15,3] = 0.2401921162199017
$arr[15,4] = 0
$arr[15,5] = 2.454531799089864
$arr[15,6] = 1.945083505773823
$arr[15,7] = 1.718452461797384
$arr[15,8] = 0.1093408658583996
$arr[15,9] = 0.6583045255649154
$arr[15,10] = 0
$arr[15,11] = 0.4196196333788151

$arr[16,0] = 0.4850129681661599
$arr[16,1] = 0.1710112026275397
$arr[16,2] = 0.639700920373599
$arr[16,3] = 0.2396498731850869
$arr[16,4] = 0
$arr[16,5] = 2.447045414080151
$arr[16,6] = 1.942664999514705
$arr[16,7] = 1.715237679479884
$arr[16,8] = 0.1091901813845553
$arr[16,9] = 0.6499948640859543
$arr[16,10] = 0
$arr[16,11] = 0.4168316980339384

$arr[17,0] = 0.4827628631303185
$arr[17,1] = 0.1704040642438258
$arr[17,2] = 0.6390866836469513
$arr[17,3] = 0.2394680026406846
$arr[17,4] = 0
$arr[17,5] = 2.4445261936836
$arr[17,6] = 1.941856490702349
$arr[17,7] = 1.71415854120832
$arr[17,8] = 0.1091398685788008
$arr[17,9] = 0.6471878538652049
$arr[17,10] = 0
$arr[17,11] = 0.4158916920895805

$arr[18,0] = 0.4929076075096077
$arr[18,1] = 0.1731431124534311
$arr[18,2] = 0.6418700853462269
$arr[18,3] = 0.2402932919145684
$arr[18,4] = 0
$arr[18,5] = 2.455924745721404
$arr[18,6] = 1.945536036527869
$arr[18,7] = 1.719051879728667
$arr[18,8] = 0.1093690896932031
$arr[18,9] = 0.6598455437456892
$arr[18,10] = 0
$arr[18,11] = 0.4201374894736034

$arr[19,0] = 0.5273249959577981
$arr[19,1] = 0.1824661756634498
$arr[19,2] = 0.6515590377996148
$arr[19,3] = 0.2431860349495523
$arr[19,4] = 0
$arr[19,5] = 2.49530524075513
$arr[19,6] = 1.958619558481786
$arr[19,7] = 1.736142814286396
$arr[19,8] = 0.1101883797091219
$arr[19,9] = 0.7028255592888115
$arr[19,10] = 0
$arr[19,11] = 0.4346768757696395

$arr[20,0] = 0.5500660012428682
$arr[20,1] = 0.1886492236913284
$arr[20,2] = 0.6581424738995736
$arr[20,3] = 0.245165939480799
$arr[20,4] = 0
$arr[20,5] = 2.521851937004641
$arr[20,6] = 1.967708001258444
$arr[20,7] = 1.747798146266959
$arr[20,8] = 0.1107604575472578
$arr[20,9] = 0.7312521438948636
$arr[20,10] = 0
$arr[20,11] = 0.444383411051291

$arr[21,0] = 0.5379061076454263
$arr[21,1] = 0.1853409721041714
$arr[21,2] = 0.6546057966221781
$arr[21,3] = 0.2441010547801667
$arr[21,4] = 0
$arr[21,5] = 2.507609279159425
$arr[21,6] = 1.96280811021478
$arr[21,7] = 1.741532997698499
$arr[21,8] = 0.1104517776054763
$arr[21,9] = 0.7160495453867384
$arr[21,10] = 0
$arr[21,11] = 0.4391841701813135

$arr[22,0] = 0.4923492546309376
$arr[22,1] = 0.1729922453658617
$arr[22,2] = 0.6417159657550258
$arr[22,3] = 0.2402475199043153
$arr[22,4] = 0
$arr[22,5] = 2.455294723160364
$arr[22,6] = 1.94533126299234
$arr[22,7] = 1.718780718141204
$arr[22,8] = 0.1093563171049396
$arr[22,9] = 0.6591487433531142
$arr[22,10] = 0
$arr[22,11] = 0.4199032992278688

$arr[23,0] = 0.4442735659873449
$arr[23,1] = 0.160056252661434
$arr[23,2] = 0.6288910743197107
$arr[23,3] = 0.2364745093104155
$arr[23,4] = 0
$arr[23,5] = 2.402325312008884
$arr[23,6] = 1.928784995967021
$arr[23,7] = 1.696316358388515
$arr[23,8] = 0.1083319180398163
$arr[23,9] = 0.599218247615056
$arr[23,10] = 0
$arr[23,11] = 0.399982918022161

$ws.Range("B2:M25").Value = $arr
